$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the existing "IP" header cell (H1) onto the
# two new header cells so they pick up the same cellXfs index (s="1")
# instead of minting new style/font entries.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for rows 2 and 3
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 5
